# Fruta / hortaliza, semanal
# Update weekly price records for Hortaliza / Haba - Terminal La Palmera de La Serena.
# The underlying data rows (Fecha, Volumen, Precio minimo/maximo/promedio, Origen,
# Precio $/Kg) are re-synced to the latest weekly extract; row order on disk stays
# the same (rows 2-14) but each row's values are refreshed to the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44376
$ws.Range("J2").Value = 400
$ws.Range("K2").Value = 12000
$ws.Range("L2").Value = 13000
$ws.Range("M2").Value = 12500
$ws.Range("O2").Value = "Provincia del Elquí"
$ws.Range("P2").Value = 500

# Row 3
$ws.Range("D3").Value = 44384
$ws.Range("J3").Value = 560
$ws.Range("K3").Value = 11500
$ws.Range("L3").Value = 12000
$ws.Range("M3").Value = 11750
$ws.Range("O3").Value = "Provincia del Elquí"
$ws.Range("P3").Value = 470

# Row 4
$ws.Range("D4").Value = 44446
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 11000
$ws.Range("L4").Value = 12000
$ws.Range("M4").Value = 11500
$ws.Range("O4").Value = "Provincia del Elquí"
$ws.Range("P4").Value = 460

# Row 5
$ws.Range("D5").Value = 44425
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 11500
$ws.Range("L5").Value = 12000
$ws.Range("M5").Value = 11750
$ws.Range("O5").Value = "Provincia del Elquí"
$ws.Range("P5").Value = 470

# Row 6
$ws.Range("D6").Value = 44377
$ws.Range("J6").Value = 520
$ws.Range("K6").Value = 12500
$ws.Range("L6").Value = 13000
$ws.Range("M6").Value = 12750
$ws.Range("O6").Value = "Provincia del Elquí"
$ws.Range("P6").Value = 510

# Row 7
$ws.Range("D7").Value = 44316
$ws.Range("J7").Value = 300
$ws.Range("K7").Value = 16000
$ws.Range("L7").Value = 17000
$ws.Range("M7").Value = 16500
$ws.Range("O7").Value = "Provincia del Elquí"
$ws.Range("P7").Value = 660

# Row 8
$ws.Range("D8").Value = 44484
$ws.Range("J8").Value = 400
$ws.Range("K8").Value = 9000
$ws.Range("L8").Value = 10000
$ws.Range("M8").Value = 9500
$ws.Range("O8").Value = "Provincia del Elquí"
$ws.Range("P8").Value = 380

# Row 9
$ws.Range("D9").Value = 44356
$ws.Range("J9").Value = 500
$ws.Range("K9").Value = 13000
$ws.Range("L9").Value = 14000
$ws.Range("M9").Value = 13500
$ws.Range("O9").Value = "Provincia de Limarí"
$ws.Range("P9").Value = 540

# Row 10
$ws.Range("D10").Value = 44370
$ws.Range("J10").Value = 520
$ws.Range("K10").Value = 13000
$ws.Range("L10").Value = 14000
$ws.Range("M10").Value = 13500
$ws.Range("O10").Value = "Provincia del Elquí"
$ws.Range("P10").Value = 540

# Row 11
$ws.Range("D11").Value = 44372
$ws.Range("J11").Value = 500
$ws.Range("K11").Value = 13000
$ws.Range("L11").Value = 14000
$ws.Range("M11").Value = 13500
$ws.Range("O11").Value = "Provincia del Elquí"
$ws.Range("P11").Value = 540

# Row 12
$ws.Range("D12").Value = 44473
$ws.Range("J12").Value = 500
$ws.Range("K12").Value = 8500
$ws.Range("L12").Value = 9000
$ws.Range("M12").Value = 8750
$ws.Range("O12").Value = "Provincia del Elquí"
$ws.Range("P12").Value = 350

# Row 13
$ws.Range("D13").Value = 44386
$ws.Range("J13").Value = 500
$ws.Range("K13").Value = 11000
$ws.Range("L13").Value = 12000
$ws.Range("M13").Value = 11500
$ws.Range("O13").Value = "Provincia del Elquí"
$ws.Range("P13").Value = 460

# Row 14
$ws.Range("D14").Value = 44466
$ws.Range("J14").Value = 400
$ws.Range("K14").Value = 9500
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = 9750
$ws.Range("O14").Value = "Provincia del Elquí"
$ws.Range("P14").Value = 390
